$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for rows 2-10, columns G-T
$data = @{
    2  = @{ G=1.017537666666667;  H=3.052613;   I=0.004254718784795717; J=0.004254718784795717;
            M=0.141694;           N=0.425082;   O=0.01763793963212447;  P=0.01763793963212447;
            Q=0.1441789821406667; R=1.297610839266; S=0.00007504447307789284; T=0.00007504447307789283 }
    3  = @{ G=1.017537666666667;  H=3.052613;   I=0.004254718784795717; J=0.004254718784795717;
            O=0.2714637835982539; P=0.2714637835982538;
            Q=2.219044447570556;  R=19.971400028135; S=0.00115500205946721; T=0.00115500205946721 }
    4  = @{ G=1.017537666666667;  H=3.052613;   I=0.004254718784795717; J=0.004254718784795717;
            M=5.710985666666667;  N=17.132957;  O=0.7108982767696218;   P=0.7108982767696217;
            Q=5.811143029626778;  R=52.30028726664101; S=0.003024672252250615; T=0.003024672252250614 }
    5  = @{ I=0.9155241810663287; J=0.9155241810663287;
            M=0.141694;           N=0.425082;   O=0.01763793963212447;  P=0.01763793963212447;
            Q=31.02422303984266;  R=279.218007358584; S=0.0161479602373981; T=0.01614796023739809 }
    6  = @{ I=0.9155241810663287; J=0.9155241810663287;
            O=0.2714637835982539; P=0.2714637835982538;
            S=0.2485316581679585; T=0.2485316581679584 }
    7  = @{ I=0.9155241810663287; J=0.9155241810663287;
            M=5.710985666666667;  N=17.132957;  O=0.7108982767696218;   P=0.7108982767696217;
            Q=1250.433279461454;  R=11253.89951515309; S=0.6508445626609723; T=0.6508445626609722 }
    8  = @{ G=19.18528466666666;  H=57.555854;  I=0.08022110014887562;  J=0.08022110014887564;
            M=0.141694;           N=0.425082;   O=0.01763793963212447;  P=0.01763793963212447;
            Q=2.718439725558666;  R=24.465957530028; S=0.00141493492164848; T=0.001414934921648479 }
    9  = @{ G=19.18528466666666;  H=57.555854;  I=0.08022110014887562;  J=0.08022110014887564;
            O=0.2714637835982539; P=0.2714637835982538;
            Q=41.83923682559222;  R=376.55313143033; S=0.02177712337082822; T=0.02177712337082822 }
    10 = @{ G=19.18528466666666;  H=57.555854;  I=0.08022110014887562;  J=0.08022110014887564;
            M=5.710985666666667;  N=17.132957;  O=0.7108982767696218;   P=0.7108982767696217;
            Q=109.5668857422531;  R=986.101971680278; S=0.05702904185639893; T=0.05702904185639893 }
}

foreach ($rowNum in $data.Keys) {
    $cols = $data[$rowNum]
    foreach ($colLetter in $cols.Keys) {
        $addr = "$colLetter$rowNum"
        $ws.Range($addr).Value = $cols[$colLetter]
    }
}
